# Updates the standard-deviation/mean/fold tables for the "RF" (row 6) and
# "Ensemble" (row 7) classifiers with newly computed dummy-data values, as
# part of adding tables for standard deviation of precision/recall and the
# PROMISE requirements dummy data set.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 6 (RF) ---
$ws.Range("B6").Value = 0.8085302072762568
$ws.Range("C6").Value = 0.0390515635286838
$ws.Range("D6").Value = 0.7569882676265656
$ws.Range("F6").Value = 0.8173729408977245
$ws.Range("G6").Value = 0.8182372478694272
$ws.Range("H6").Value = 0.8711342581896228
$ws.Range("I6").Value = 0.8256457789346012
$ws.Range("J6").Value = 0.01209272630147904
$ws.Range("K6").Value = 0.8216235142705731
$ws.Range("L6").Value = 0.8209096459096459
$ws.Range("M6").Value = 0.8180666720564853
$ws.Range("N6").Value = 0.8179769068832798
$ws.Range("O6").Value = 0.849652155553021
$ws.Range("P6").Value = 0.7513773768325669
$ws.Range("Q6").Value = 0.05310457069409973
$ws.Range("R6").Value = 0.7290270127950293
$ws.Range("S6").Value = 0.6763706634674377
$ws.Range("T6").Value = 0.827485368004644
$ws.Range("U6").Value = 0.7307229880007221
$ws.Range("V6").Value = 0.7932808518950014
$ws.Range("W6").Value = 0.8095402314725835
$ws.Range("X6").Value = 0.01323351489720702
$ws.Range("Y6").Value = 0.7881597495340692
$ws.Range("Z6").Value = 0.8070828896035197
$ws.Range("AA6").Value = 0.8068226216452928
$ws.Range("AB6").Value = 0.8176951992844992
$ws.Range("AC6").Value = 0.827940697295536
$ws.Range("AD6").Value = 0.8030358157294298
$ws.Range("AE6").Value = 0.02877332982882458
$ws.Range("AF6").Value = 0.7701091461729759
$ws.Range("AG6").Value = 0.7779402669046132
$ws.Range("AH6").Value = 0.7992824395625754
$ws.Range("AI6").Value = 0.8180590213800699
$ws.Range("AJ6").Value = 0.8497882046269143
$ws.Range("AK6").Value = 0.8189970648453558
$ws.Range("AL6").Value = 0.046231031292413
$ws.Range("AM6").Value = 0.7837671870105706
$ws.Range("AN6").Value = 0.7642065819485175
$ws.Range("AO6").Value = 0.8491551294018238
$ws.Range("AP6").Value = 0.8056723240043461
$ws.Range("AQ6").Value = 0.8921841018615212

# --- Row 7 (Ensemble) ---
$ws.Range("B7").Value = 0.8460971644411467
$ws.Range("C7").Value = 0.05013449088806315
$ws.Range("D7").Value = 0.8410667968338742
$ws.Range("E7").Value = 0.7662457743102904
$ws.Range("G7").Value = 0.828003280576886
$ws.Range("H7").Value = 0.9139619013733883
$ws.Range("I7").Value = 0.8618252690317654
$ws.Range("J7").Value = 0.03216494561268367
$ws.Range("M7").Value = 0.8696912686725929
$ws.Range("P7").Value = 0.8313807466572986
$ws.Range("Q7").Value = 0.04099386412136994
$ws.Range("R7").Value = 0.7996616223469558
$ws.Range("S7").Value = 0.808845414108572
$ws.Range("T7").Value = 0.8513269831642902
$ws.Range("V7").Value = 0.9028587055337645
$ws.Range("W7").Value = 0.8499449907477221
$ws.Range("X7").Value = 0.03733860504644549
$ws.Range("Y7").Value = 0.8401035523894891
$ws.Range("Z7").Value = 0.7862047714040123
$ws.Range("AC7").Value = 0.891832719252074
$ws.Range("AD7").Value = 0.8689077229551714
$ws.Range("AE7").Value = 0.04264945786907141
$ws.Range("AF7").Value = 0.8940831774608371
$ws.Range("AG7").Value = 0.7952772506852204
$ws.Range("AH7").Value = 0.893185533104888
$ws.Range("AI7").Value = 0.8481485675034062
$ws.Range("AJ7").Value = 0.9138440860215054
$ws.Range("AK7").Value = 0.8351239298215546
$ws.Range("AL7").Value = 0.03103274933424498
$ws.Range("AM7").Value = 0.8411041475390205
$ws.Range("AN7").Value = 0.7855013878782979
$ws.Range("AO7").Value = 0.8404567206986562
$ws.Range("AP7").Value = 0.8263972586553232
$ws.Range("AQ7").Value = 0.882160134336475
